$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# ALC
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376

$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880

$ws.Range("H97").Value = 666
$ws.Range("I97").Value = 499
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1497
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1001
$ws.Range("N97").Value = -3992

$ws.Range("H112").Value = 922.8
$ws.Range("J112").Value = 922.8
$ws.Range("L112").Value = 2768.4
$ws.Range("N112").Value = -4984.4

$ws.Range("H116").Value = 5000
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884

$ws.Range("H121").Value = 1000
$ws.Range("J121").Value = 1000
$ws.Range("L121").Value = 3000
$ws.Range("N121").Value = -6494

$ws.Range("H138").Value = 5126.7886
$ws.Range("J138").Value = 6128.3657
$ws.Range("L138").Value = 18385.0971
$ws.Range("N138").Value = -28665.0971

# -----------------------------------------------------------------
# ARM
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H124").Value = 33370034
$ws.Range("J124").Value = 33370034
$ws.Range("L124").Value = 33370034
$ws.Range("N124").Value = -33379854

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140

# -----------------------------------------------------------------
# BSM
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H94").Value = 965.1111
$ws.Range("I94").Value = 781.1667
$ws.Range("K94").Value = 781.1667
$ws.Range("M94").Value = -330.1667

# Rows 117-141 all lose their currentAveragePrice/Leve price/profit
# columns (H:N), except rows 121 and 136 which stay untouched.
$skipRows = @(121, 136)
for ($r = 117; $r -le 141; $r++) {
    if ($skipRows -contains $r) { continue }
    $ws.Range("H$r`:N$r").ClearContents()
}

# -----------------------------------------------------------------
# CRP
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H28").Value = 100643
$ws.Range("J28").Value = 100643
$ws.Range("L28").Value = 100643
$ws.Range("N28").Value = -101133

$ws.Range("H103").Value = 36663
$ws.Range("I103").Value = 29994.5
$ws.Range("J103").Value = 50000
$ws.Range("K103").Value = 29994.5
$ws.Range("L103").Value = 50000
$ws.Range("M103").Value = -28822.5
$ws.Range("N103").Value = -52344

# -----------------------------------------------------------------
# CUL
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 575.8570999999999
$ws.Range("J12").Value = 511.75
$ws.Range("L12").Value = 1535.25
$ws.Range("N12").Value = -1881.25

# -----------------------------------------------------------------
# GSM
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 5206.8335
$ws.Range("I70").Value = 5206.8335
$ws.Range("K70").Value = 5206.8335
$ws.Range("M70").Value = -4936.8335

$ws.Range("H73").Value = 5206.8335
$ws.Range("I73").Value = 5206.8335
$ws.Range("K73").Value = 5206.8335
$ws.Range("M73").Value = -4270.8335

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H132").Value = 5996.6
$ws.Range("I132").Value = 5996.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17989.8
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15459.8
$ws.Range("N132").ClearContents()

# -----------------------------------------------------------------
# WVR
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 9875
$ws.Range("I62").Value = 8500
$ws.Range("J62").Value = 10333.333
$ws.Range("K62").Value = 8500
$ws.Range("L62").Value = 10333.333
$ws.Range("M62").Value = -7876
$ws.Range("N62").Value = -11581.333

$ws.Range("H65").Value = 9875
$ws.Range("I65").Value = 8500
$ws.Range("J65").Value = 10333.333
$ws.Range("K65").Value = 42500
$ws.Range("L65").Value = 51666.665
$ws.Range("M65").Value = -39380
$ws.Range("N65").Value = -57906.665
